# Appends 16 new trading-log rows (60-75) to Sheet1, mirroring the
# Python trading bot run logged at 2025-10-04T01:22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("timestamp","action","token","signal_type","price","position_size_usd","leverage","stiffness","pnl_percent","exit_reason","status","error_message")

$rows = @(
    @(60, "2025-10-04T01:22:14.291475", "TRADING_ATTEMPT", "SOL", "UNKNOWN", 231.386880157209, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 1/8"),
    @(61, "2025-10-04T01:22:16.169141", "POSITION_FAILED", "SOL", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 1"),
    @(62, "2025-10-04T01:22:16.196986", "TRADING_ATTEMPT", "BTC", "UNKNOWN", 122145.0094321105, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 2/8"),
    @(63, "2025-10-04T01:22:17.935069", "POSITION_OPENED", "BTC", "UNKNOWN", 122145.0094321105, 90, 1, 0, $null, $null, "SUCCESS", $null),
    @(64, "2025-10-04T01:22:17.961778", "TRADING_ATTEMPT", "ETH", "UNKNOWN", 4492.563675941279, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 3/8"),
    @(65, "2025-10-04T01:22:19.648871", "POSITION_FAILED", "ETH", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 3"),
    @(66, "2025-10-04T01:22:19.677577", "TRADING_ATTEMPT", "ARB", "UNKNOWN", 0.4475369841415226, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 4/8"),
    @(67, "2025-10-04T01:22:21.331050", "POSITION_FAILED", "ARB", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 4"),
    @(68, "2025-10-04T01:22:21.361777", "TRADING_ATTEMPT", "SUI", "UNKNOWN", 3.580526691599361, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 5/8"),
    @(69, "2025-10-04T01:22:23.144754", "POSITION_FAILED", "SUI", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 5"),
    @(70, "2025-10-04T01:22:23.174530", "TRADING_ATTEMPT", "XRP", "UNKNOWN", 3.037337837354311, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 6/8"),
    @(71, "2025-10-04T01:22:24.790784", "POSITION_FAILED", "XRP", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 6"),
    @(72, "2025-10-04T01:22:24.820697", "TRADING_ATTEMPT", "AAVE", "UNKNOWN", 289.7331155048397, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 7/8"),
    @(73, "2025-10-04T01:22:26.497201", "POSITION_FAILED", "AAVE", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 7"),
    @(74, "2025-10-04T01:22:26.526928", "TRADING_ATTEMPT", "ENA", "UNKNOWN", 0.6116338841784414, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 8/8"),
    @(75, "2025-10-04T01:22:28.293407", "POSITION_FAILED", "ENA", "UNKNOWN", $null, $null, $null, $null, $null, $null, "FAILED", "Trade execution failed for trade 8")
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $val = $row[$i + 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $i + 1).Value = $val
        }
    }
}
